$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 11 (old) -> data logically shifts: old row10 content moves to row9,
#    old row11's "C" entry moves into row10's C cell, and rows 7/8 gain new
#    bitácora entries. Row 11 itself ends up with only the index cell (A11).
# ---------------------------------------------------------------------------

# Completely clear B11:F11 (content + formatting) so no <c> nodes remain
# there except A11, matching the target (row 11 keeps only its index value).
$ws.Range("B11:F11").Clear()

# ---------------------------------------------------------------------------
# 2) Row 10: keep A10/B10 as-is; change C10's text, bump D10's date, and give
#    E10/F10 new text (these two are first-time strings -> allocate shared
#    string indices 15 and 16, so we write them before any other new text).
# ---------------------------------------------------------------------------

# E10/F10 currently hold border-only "left/top wrap+shrink" formatting (s=9)
# already -- just change their text (first brand-new strings written).
$ws.Range("E10").Value = "Se observó que la función no tenía ningún parámetro que pusiera la información dentro del HTML, por lo que se colocarón parámetros que buscaran los datos en la lista infoPaises."
$ws.Range("F10").Value = "2 días"

# F10 used to carry the bold/font variant of the "left/top wrap+shrink" style;
# re-point it at the plain variant (same alignment, no extra font) by copying
# formats from an existing cell that already has that exact style (E6).
$ws.Range("E6").Copy()
$ws.Range("F10").PasteSpecial(-4122) | Out-Null
$ws.Range("F10").Value = "2 días"

# C10 switches text to the old C11 content and adopts the "left/top
# wrap+shrink" style (same as E10/E6/F6 etc.).
$ws.Range("F6").Copy()
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value = "En la función hablantes por idioma genera el HTML vacío."

# D10 stays on the date style it already has; only the date value changes.
$ws.Range("D10").Value = 45203

# ---------------------------------------------------------------------------
# 3) Row 9: becomes what row 10 used to contain (C/E/F text + B/D dates).
# ---------------------------------------------------------------------------

# B9 / D9 need the date style (numFmt 14, left/top, wrap+shrink) - copy it
# from B6 which already carries that exact style.
$ws.Range("B6").Copy()
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = 45201

$ws.Range("B6").Copy()
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = 45201

# C9/E9 keep the plain border/wrap/shrink style they already have (s=3).
$ws.Range("C9").Value = "En la función de mostrar los códigos de un país, da error con cualquier selección del país"
$ws.Range("E9").Value = 'Se eliminó el argumento "inicio=1" de las llamadas a la función "enumerate" y se ajustaron los índices dentro del bucle "for" para que comiencen desde 1.'

# F9 adopts the "left/top wrap+shrink" style (copy from F6) and the "1.5
# Horas" text.
$ws.Range("F6").Copy()
$ws.Range("F9").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Value = "1.5 Horas"

$ws.Rows(9).RowHeight = 47.25

# ---------------------------------------------------------------------------
# 4) Row 7: brand-new bitácora entry (second batch of new strings -> indices
#    17,18,19).
# ---------------------------------------------------------------------------

$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Value = 45194

$ws.Range("C7").Value = "Error al generar un HTML"

# D7 needs a NEW style: same plain border/wrap/shrink alignment it already
# has (s=3) plus a date number format - applying NumberFormat keeps the
# existing alignment and creates (or reuses) exactly that combination.
$ws.Range("D7").Value = 45194
$ws.Range("D7").NumberFormat = "mm-dd-yy"

$ws.Range("F6").Copy()
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = "Se observó que era más fácil trabajar con la librería BeautifulSoup4, con la finalidad de crear archivos HTML más fácilmente.  "

$ws.Range("F7").Value = "3 horas"

$ws.Rows(7).RowHeight = 48.75

# ---------------------------------------------------------------------------
# 5) Row 8: brand-new bitácora entry (final new strings -> indices 20,21).
# ---------------------------------------------------------------------------

$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = 45200

$ws.Range("C8").Value = "El HTML no se abre al corres cualquier función"

$ws.Range("D8").Value = 45203
$ws.Range("D8").NumberFormat = "mm-dd-yy"

$ws.Range("F6").Copy()
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = "Se observó que se debía importar la librería webbrowser para que el return abra el HTML en el navegador de preferencia del usuario"

# F8 stays blank (already s=3, no content) - nothing to do.

$ws.Rows(8).RowHeight = 56.25

# ---------------------------------------------------------------------------
# 6) Row 11 final height.
# ---------------------------------------------------------------------------
$ws.Rows(11).RowHeight = 54.75

# ---------------------------------------------------------------------------
# 7) Column D a touch wider, drop the "best fit" flag.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 12.86

# ---------------------------------------------------------------------------
# 8) Selection / view: no more frozen top-left scroll, select J6 instead of
#    F11.
# ---------------------------------------------------------------------------
$ws.Range("J6").Select()
